# Thêm cột "Địa chỉ thửa đất" vào báo cáo quyết định miễn tiền thuê đất.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Chèn một cột mới trước cột D (Excel sẽ tự đẩy các cột D:K sang E:L
# và sao chép định dạng/merge theo kiểu "Insert Entire Column" chuẩn).
$ws.Columns("D").Insert() | Out-Null

# Cột mới lấy độ rộng giống cột C liền kề.
$ws.Columns("D").ColumnWidth = $ws.Columns("C").ColumnWidth

# Đặt tiêu đề cho cột mới và gộp ô tiêu đề (D5:D6) giống các cột tiêu đề khác.
$ws.Range("D5").Value = "Địa chỉ thửa đất"
$ws.Range("D5:D6").Merge() | Out-Null

# Cập nhật vùng chọn hiện hành như trong bản gốc.
$ws.Range("E3").Select() | Out-Null
